$wb = $excel.ActiveWorkbook

# New sheet names, in the same order as the existing sheets (position-based rename).
$newNames = @(
    "summ12907143",
    "summ13109003",
    "summ13364027",
    "summ13640579",
    "summ13913182",
    "summ14172226",
    "summ14416751",
    "summ14671432",
    "summ14923992",
    "summ15192799",
    "summ15455335",
    "summ15734916",
    "summ16049363",
    "summ16478982",
    "summ16750526",
    "summ17017297",
    "summ17311384",
    "summ17657499",
    "summ17930693",
    "summ18191849",
    "summ18453832",
    "summ18691225",
    "summ18935997",
    "summ19195541",
    "summ19464076",
    "summ19718891",
    "summ19972220",
    "summ20220757",
    "summ20464910",
    "summ20723651"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
